$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A104").Value = 1.963
$ws.Range("B104").Value = 1.578
$ws.Range("C104").Value = 1.699

$ws.Range("A105").Value = 2.218
$ws.Range("B105").Value = 1.909
$ws.Range("C105").Value = 1.949
